$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10. This shifts the existing rows 10-29 down to 11-30,
# carrying their formatting (including the date style on column D) along with them.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with this week's new record.
$ws.Cells.Item(10, 1).Value = 4
$ws.Cells.Item(10, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(10, 3).Value = "Los Lagos"
$ws.Cells.Item(10, 4).Value = 44757
$ws.Cells.Item(10, 5).Value = 10
$ws.Cells.Item(10, 6).Value = 100112012
$ws.Cells.Item(10, 7).Value = "Espinaca"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 35
$ws.Cells.Item(10, 11).Value = 13000
$ws.Cells.Item(10, 12).Value = 13000
$ws.Cells.Item(10, 13).Value = 13000
$ws.Cells.Item(10, 14).Value = '$/cuna 10 kilos'
$ws.Cells.Item(10, 15).Value = "Región Metropolitana"
$ws.Cells.Item(10, 16).Value = 1300
$ws.Cells.Item(10, 17).Value = 10
$ws.Cells.Item(10, 18).Value = "Hortaliza"
